# Excel COM-interop script implementing commit "Assets and Tracking Script added".
#
# The loading-time tracking sheet gains seven new tracked days (columns W:BS),
# each day occupying a variable number of columns (one per asset/script tracked
# that day), mirroring the existing layout in columns B:V. Row 1 holds the date
# label for each day; rows 2-7 hold the per-asset metric values for that day.
#
# Because column A:XFD in this sheet defaults to the "Text" number format, newly
# written numeric cells must first be reset to the "Normal" style so that Excel
# stores them as real numbers (t="n") instead of text - exactly like the original
# numeric cells in columns B:V.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.Style = "Normal"
    $cell.Value = $value
}

# Row 1: date headers for the new days ( Oct 26,  Oct 27,  Oct 28,  Nov 01,  Nov 02,  Nov 03,  Nov 07)
$row1 = @(
    @("W1", " Oct 26"),
    @("X1", " Oct 26"),
    @("Y1", " Oct 26"),
    @("Z1", " Oct 26"),
    @("AA1", " Oct 27"),
    @("AB1", " Oct 28"),
    @("AC1", " Oct 28"),
    @("AD1", " Oct 28"),
    @("AE1", " Oct 28"),
    @("AF1", " Oct 28"),
    @("AG1", " Oct 28"),
    @("AH1", " Oct 28"),
    @("AI1", " Nov 01"),
    @("AJ1", " Nov 01"),
    @("AK1", " Nov 01"),
    @("AL1", " Nov 01"),
    @("AM1", " Nov 01"),
    @("AN1", " Nov 01"),
    @("AO1", " Nov 01"),
    @("AP1", " Nov 01"),
    @("AQ1", " Nov 01"),
    @("AR1", " Nov 01"),
    @("AS1", " Nov 02"),
    @("AT1", " Nov 02"),
    @("AU1", " Nov 02"),
    @("AV1", " Nov 02"),
    @("AW1", " Nov 02"),
    @("AX1", " Nov 03"),
    @("AY1", " Nov 03"),
    @("AZ1", " Nov 03"),
    @("BA1", " Nov 03"),
    @("BB1", " Nov 03"),
    @("BC1", " Nov 03"),
    @("BD1", " Nov 03"),
    @("BE1", " Nov 03"),
    @("BF1", " Nov 03"),
    @("BG1", " Nov 07"),
    @("BH1", " Nov 07"),
    @("BI1", " Nov 07"),
    @("BJ1", " Nov 07"),
    @("BK1", " Nov 07"),
    @("BL1", " Nov 07"),
    @("BM1", " Nov 07"),
    @("BN1", " Nov 07"),
    @("BO1", " Nov 07"),
    @("BP1", " Nov 07"),
    @("BQ1", " Nov 07"),
    @("BR1", " Nov 07"),
    @("BS1", " Nov 07")
)
foreach ($pair in $row1) { Set-Cell $pair[0] $pair[1] }

# Row 2: Login metric for each new day
$row2 = @(
    @("W2", 0),
    @("X2", 0),
    @("Y2", 0),
    @("Z2", 0),
    @("AA2", 0),
    @("AB2", 0),
    @("AC2", 0),
    @("AD2", 0),
    @("AE2", 0),
    @("AF2", 0),
    @("AG2", 0),
    @("AH2", 0),
    @("AO2", 0),
    @("AR2", 0),
    @("AS2", 0),
    @("AT2", 0),
    @("AU2", 0),
    @("AV2", 0),
    @("AW2", 0),
    @("BI2", 0),
    @("BJ2", 0),
    @("BK2", 0),
    @("BL2", 0),
    @("BM2", 0)
)
foreach ($pair in $row2) { Set-Cell $pair[0] $pair[1] }

# Row 3: Dashboard metric for each new day
$row3 = @(
    @("W3", 7),
    @("X3", 8),
    @("Y3", 8),
    @("Z3", 8),
    @("AA3", 8),
    @("AB3", 8),
    @("AC3", 7),
    @("AD3", 10),
    @("AE3", 8),
    @("AF3", 10),
    @("AG3", 10),
    @("AH3", 8),
    @("AJ3", 10),
    @("AL3", 10),
    @("AN3", 11),
    @("AO3", 9),
    @("AP3", 11),
    @("AQ3", 12),
    @("AR3", 10),
    @("AS3", 9),
    @("AT3", 10),
    @("AU3", 8),
    @("AV3", 10),
    @("AW3", 10),
    @("AX3", 36),
    @("AY3", 36),
    @("AZ3", 36),
    @("BA3", 36),
    @("BB3", 36),
    @("BC3", 10),
    @("BD3", 10),
    @("BE3", 7),
    @("BF3", 8),
    @("BG3", 9),
    @("BH3", 9),
    @("BI3", 10),
    @("BJ3", 9),
    @("BK3", 9),
    @("BL3", 9),
    @("BM3", 9),
    @("BN3", 9),
    @("BO3", 9),
    @("BP3", 10),
    @("BQ3", 10),
    @("BR3", 9),
    @("BS3", 9)
)
foreach ($pair in $row3) { Set-Cell $pair[0] $pair[1] }

# Row 4: Trucks metric for each new day
$row4 = @(
    @("W4", 19),
    @("X4", 36),
    @("Y4", 9),
    @("Z4", 7),
    @("AA4", 8),
    @("AB4", 23),
    @("AC4", 71),
    @("AD4", 31),
    @("AE4", 58),
    @("AF4", 25),
    @("AG4", 258),
    @("AH4", 19),
    @("AI4", 57),
    @("AJ4", 45),
    @("AK4", 11),
    @("AL4", 24),
    @("AM4", 26),
    @("AN4", 20),
    @("AO4", 21),
    @("AP4", 32),
    @("AQ4", 21),
    @("AR4", 53),
    @("AS4", 65),
    @("AT4", 42),
    @("AU4", 30),
    @("AV4", 47),
    @("AW4", 93),
    @("AX4", 11),
    @("AY4", 22),
    @("AZ4", 60),
    @("BA4", 24),
    @("BB4", 18),
    @("BC4", 25),
    @("BD4", 99),
    @("BE4", 107),
    @("BF4", 177),
    @("BG4", 25),
    @("BH4", 24),
    @("BI4", 29),
    @("BJ4", 32),
    @("BK4", 25),
    @("BL4", 26),
    @("BM4", 25),
    @("BN4", 25),
    @("BO4", 26),
    @("BP4", 25),
    @("BQ4", 26),
    @("BR4", 28),
    @("BS4", 25)
)
foreach ($pair in $row4) { Set-Cell $pair[0] $pair[1] }

# Row 5: Deleted Trailers metric for each new day
$row5 = @(
    @("Z5", 0),
    @("AA5", 0),
    @("AB5", 0),
    @("AC5", 0),
    @("AD5", 0),
    @("AE5", 0),
    @("AF5", 0),
    @("AG5", 0),
    @("AH5", 0),
    @("AO5", 0),
    @("AR5", 0),
    @("AS5", 0),
    @("AT5", 0),
    @("AU5", 0),
    @("AV5", 0),
    @("AW5", 0),
    @("BI5", 0),
    @("BJ5", 0),
    @("BK5", 0),
    @("BL5", 0),
    @("BM5", 0)
)
foreach ($pair in $row5) { Set-Cell $pair[0] $pair[1] }

# Row 6: Deleted Trucks metric for each new day
$row6 = @(
    @("Z6", 0),
    @("AA6", 0),
    @("AB6", 0),
    @("AD6", 0),
    @("AE6", 0),
    @("AF6", 0),
    @("AG6", 7),
    @("AH6", 0),
    @("AO6", 0),
    @("AR6", 0),
    @("AS6", 0),
    @("AT6", 1),
    @("AU6", 0),
    @("AV6", 0),
    @("AW6", 0),
    @("BI6", 0),
    @("BJ6", 0),
    @("BK6", 0),
    @("BL6", 0),
    @("BM6", 0)
)
foreach ($pair in $row6) { Set-Cell $pair[0] $pair[1] }

# Row 7: Trailer metric for each new day
$row7 = @(
    @("Z7", 0),
    @("AA7", 0),
    @("AB7", 0),
    @("AC7", 0),
    @("AD7", 0),
    @("AE7", 0),
    @("AF7", 0),
    @("AG7", 0),
    @("AH7", 0),
    @("AN7", 0),
    @("AO7", 0),
    @("AR7", 0),
    @("AS7", 0),
    @("AT7", 0),
    @("AU7", 0),
    @("AV7", 0),
    @("AW7", 0),
    @("BI7", 0),
    @("BJ7", 0),
    @("BK7", 0),
    @("BL7", 0),
    @("BM7", 0)
)
foreach ($pair in $row7) { Set-Cell $pair[0] $pair[1] }
